$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Overview")

# Row 27: "سود هر سهم بر اساس آخرین سرمایه" (EPS based on latest capital)
# Updated per the new read_price algorithm.
$ws.Range("D27").Value = 162
$ws.Range("E27").Value = 290
$ws.Range("F27").Value = 434
$ws.Range("G27").Value = 158
$ws.Range("H27").Value = 405
$ws.Range("I27").Value = 630
$ws.Range("J27").Value = 834
$ws.Range("K27").Value = 205
$ws.Range("L27").Value = 370
$ws.Range("M27").Value = 517
